# Remove spaces and quotation marks from the header cells of Sheet1.
# Rows 1-3 hold two-line (sometimes three-line) column headers built out of
# shared strings that were padded with trailing/leading spaces; row 4 holds
# the dashed separator line. This trims the stray whitespace and merges the
# "VE/" + "VO2"/"VCO2" split labels into single "VE/VO2" / "VE/CO2" headers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: main header labels -------------------------------------------------
# (M1/N1 are fixed up last, below, after the rest of the sheet has been
# de-whitespaced -- that mirrors how the edit was actually made: the split
# "VE/" + "VO2"/"VCO2" labels were merged into single cells as a separate,
# final pass.)
$ws.Range("A1").Value = "TIME"
$ws.Range("B1").Value = "VO2"
$ws.Range("C1").Value = "VO2/kg"
$ws.Range("D1").Value = "METS"
$ws.Range("E1").Value = "VCO2"
$ws.Range("F1").Value = "VE"
$ws.Range("G1").Value = "RER"
$ws.Range("H1").Value = "RR"
$ws.Range("I1").Value = "Vt"
$ws.Range("J1").Value = "FEO2"
$ws.Range("K1").Value = "FECO2"
$ws.Range("L1").Value = "HR"
$ws.Range("O1").Value = "PetCO2"
$ws.Range("P1").Value = "Power"

# --- Row 2: secondary header labels (mostly blank now) -------------------------
$ws.Range("A2").Value = ""
$ws.Range("B2").Value = "STPD"
$ws.Range("C2").Value = "STPD"
$ws.Range("D2").Value = ""
$ws.Range("E2").Value = "STPD"
$ws.Range("F2").Value = "BTPS"
$ws.Range("G2").Value = ""
$ws.Range("H2").Value = ""
$ws.Range("I2").Value = "BTPS"
$ws.Range("J2").Value = ""
$ws.Range("K2").Value = ""
$ws.Range("L2").Value = ""
$ws.Range("M2").Value = ""
$ws.Range("N2").Value = ""
$ws.Range("O2").Value = ""

# --- Row 3: unit labels ---------------------------------------------------------
$ws.Range("A3").Value = "min"
$ws.Range("B3").Value = "L/min"
$ws.Range("C3").Value = "ml/kg/m"
$ws.Range("D3").Value = ""
$ws.Range("E3").Value = "L/min"
$ws.Range("F3").Value = "L/min"
$ws.Range("G3").Value = ""
$ws.Range("H3").Value = "BPM"
$ws.Range("I3").Value = "L"
$ws.Range("J3").Value = "%"
$ws.Range("K3").Value = "%"
$ws.Range("L3").Value = "bpm"
$ws.Range("M3").Value = "BT/ST"
$ws.Range("N3").Value = "BT/ST"
$ws.Range("O3").Value = "mmHg"
$ws.Range("P3").Value = "Watts"

# --- Row 4: separator row --------------------------------------------------------
$ws.Range("A4").Value = "----------"

# --- Row 1 (cont'd): merge the split "VE/" + "VO2"/"VCO2" labels ---------------
$ws.Range("N1").Value = "VE/CO2"
$ws.Range("M1").Value = "VE/VO2"

# Update the active selection to match the saved view state.
$ws.Range("M2:N2").Select()
